$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (date serial 44706) is inserted at the top of this
# market's data block (rows 331-332), pushing every existing row down by
# two rows (the block runs through row 352, which becomes 354).
$ws.Rows("331:332").Insert()

# Row 331 - "Primera" quality for the new date.
$ws.Cells.Item(331, 1).Value = 1
$ws.Cells.Item(331, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(331, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(331, 4).Value = 44706
$ws.Cells.Item(331, 5).Value = 15
$ws.Cells.Item(331, 6).Value = 100112032
$ws.Cells.Item(331, 7).Value = "Zapallo italiano"
$ws.Cells.Item(331, 8).Value = "Huracán"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 130
$ws.Cells.Item(331, 11).Value = 9000
$ws.Cells.Item(331, 12).Value = 10000
$ws.Cells.Item(331, 13).Value = 9500
$ws.Cells.Item(331, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(331, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(331, 16).Value = 136
$ws.Cells.Item(331, 17).Value = 70
$ws.Cells.Item(331, 18).Value = "Hortaliza"

# Row 332 - "Segunda" quality for the new date.
$ws.Cells.Item(332, 1).Value = 1
$ws.Cells.Item(332, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(332, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(332, 4).Value = 44706
$ws.Cells.Item(332, 5).Value = 15
$ws.Cells.Item(332, 6).Value = 100112032
$ws.Cells.Item(332, 7).Value = "Zapallo italiano"
$ws.Cells.Item(332, 8).Value = "Huracán"
$ws.Cells.Item(332, 9).Value = "Segunda"
$ws.Cells.Item(332, 10).Value = 140
$ws.Cells.Item(332, 11).Value = 7000
$ws.Cells.Item(332, 12).Value = 8000
$ws.Cells.Item(332, 13).Value = 7500
$ws.Cells.Item(332, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(332, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(332, 16).Value = 75
$ws.Cells.Item(332, 17).Value = 100
$ws.Cells.Item(332, 18).Value = "Hortaliza"
